# "fixed monte carlo for ties"
# Updates the Monte Carlo simulation results on the "Playoff Odds" sheet
# (place-finish probability distribution + playoff chance) and the
# corresponding "Expected_Final_Record" / "Playoff_Chance_Pct" values on
# the "Record Odds" sheet.

$wb = $excel.ActiveWorkbook

$wsPlayoff = $wb.Worksheets.Item("Playoff Odds")
$wsRecord  = $wb.Worksheets.Item("Record Odds")

# ---------------------------------------------------------------------
# "Playoff Odds" sheet - columns B..L for rows 2..11
# B..K = chance (%) of finishing in that place (1st..10th)
# L    = overall chance (%) of making the playoffs
# ---------------------------------------------------------------------

$playoffData = @{
    2  = @(45.1, 28,   15.6, 7.2,  3,    0.8,  0.2,  0.1,  0,    0,    99.7)
    3  = @(30,   30.7, 22.6, 9,    4.1,  2.4,  0.8,  0.1,  0.3,  0,    98.8)
    4  = @(16.9, 22.2, 21.8, 20.2, 11.2, 3.8,  2.4,  1.1,  0.3,  0.1,  96.09999999999999)
    5  = @(6,    11.8, 20.3, 24.8, 18.7, 8.699999999999999, 5.1, 2.9, 1.2, 0.5, 90.3)
    6  = @(1.5,  4.2,  10.1, 19,   19.2, 17.7, 14.9, 8,    3.5,  1.9,  71.7)
    7  = @(0.2,  1.9,  3.8,  7.6,  15.6, 24.1, 17.3, 13.2, 9.9,  6.4,  53.2)
    8  = @(0.1,  0.2,  1.7,  4,    10.7, 15.1, 18.5, 21.2, 16.9, 11.6, 31.8)
    9  = @(0.2,  0.7,  2.2,  4.6,  8.6,  10.2, 14.5, 16.3, 22.6, 20.1, 26.5)
    10 = @(0,    0.1,  1.1,  2.1,  5.5,  8.699999999999999, 13.5, 18.2, 20.5, 30.3, 17.5)
    11 = @(0,    0.2,  0.8,  1.5,  3.4,  8.5,  12.8, 18.9, 24.8, 29.1, 14.4)
}

foreach ($row in 2..11) {
    $vals = $playoffData[$row]
    # columns B(2) through L(12)
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $col = 2 + $i
        $wsPlayoff.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------------
# "Record Odds" sheet - column G (Expected_Final_Record) for rows 2..11
# and column F (Playoff_Chance_Pct) for rows 4..11 (2 & 3 unchanged)
# ---------------------------------------------------------------------

$expectedRecord = @{
    2  = "10.3-3.6-0.1"
    3  = "9.5-4.4-0.1"
    4  = "8.8-5.2-0.1"
    5  = "7.8-6.1-0.1"
    6  = "7.0-6.9-0.1"
    7  = "6.0-7.9-0.1"
    8  = "5.5-8.3-0.1"
    9  = "5.2-8.7"
    10 = "4.8-9.2-0.1"
    11 = "4.7-9.3-0.1"
}

foreach ($row in 2..11) {
    $wsRecord.Cells.Item($row, 7).Value = $expectedRecord[$row]
}

$playoffChancePct = @{
    4  = 96.09999999999999
    5  = 90.3
    6  = 71.7
    7  = 53.2
    8  = 31.8
    9  = 26.5
    10 = 17.5
    11 = 14.4
}

foreach ($row in 4..11) {
    $wsRecord.Cells.Item($row, 6).Value = $playoffChancePct[$row]
}
